# Updates cryptos list values (Price column D, Volume(1h) column E)
# per upstream scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.081.70"
$ws.Range("E2").Value = "  +1.72%  "

$ws.Range("D3").Value = "1.959.82"
$ws.Range("E3").Value = "  -0.09%  "

$ws.Range("E4").Value = "  +0.50%  "

$ws.Range("D5").Value = "'244.94"
$ws.Range("E5").Value = "  -1.04%  "

$ws.Range("E6").Value = "  +0.33%  "

$ws.Range("D7").Value = "'0.4891"
$ws.Range("E7").Value = "  +1.75%  "

$ws.Range("D8").Value = "'0.2956"
$ws.Range("E8").Value = "  +1.49%  "

$ws.Range("D9").Value = "'0.07018"
$ws.Range("E9").Value = "  +4.23%  "

$ws.Range("D10").Value = "'19.43"
$ws.Range("E10").Value = "  +1.91%  "

$ws.Range("D11").Value = "'107.80"
$ws.Range("E11").Value = "  -0.87%  "

$ws.Range("D12").Value = "1.959.98"
$ws.Range("E12").Value = "  -0.06%  "

$ws.Range("D13").Value = "'0.07788"
$ws.Range("E13").Value = "  +1.08%  "

$ws.Range("D14").Value = "'5.486"
$ws.Range("E14").Value = "  +1.76%  "

$ws.Range("D15").Value = "'0.7011"
$ws.Range("E15").Value = "  +1.49%  "

$ws.Range("D16").Value = "'283.02"
$ws.Range("E16").Value = "  -1.83%  "

$ws.Range("D17").Value = "31.089.17"
$ws.Range("E17").Value = "  +1.65%  "

$ws.Range("D18").Value = "'13.29"
$ws.Range("E18").Value = "  +1.32%  "

$ws.Range("D19").Value = "'0.000007776"
$ws.Range("E19").Value = "  +1.75%  "

$ws.Range("D20").Value = "2.220.10"
$ws.Range("E20").Value = "  +0.27%  "

$ws.Range("D22").Value = "'5.531"
$ws.Range("E22").Value = "  -1.50%  "

$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.51%  "

$ws.Range("E25").Value = "  -0.07%  "

$ws.Range("D26").Value = "'168.39"
$ws.Range("E26").Value = "  -1.50%  "

$ws.Range("D27").Value = "'20.01"
$ws.Range("E27").Value = "  +0.67%  "

$ws.Range("D28").Value = "'2.200"
$ws.Range("E28").Value = "  +1.29%  "

$ws.Range("D29").Value = "'0.1051"
$ws.Range("E29").Value = "  -1.09%  "

$ws.Range("E30").Value = "  -3.03%  "

$ws.Range("D31").Value = "'1.580"
$ws.Range("E31").Value = "  -0.59%  "

$ws.Range("D32").Value = "'4.623"
$ws.Range("E32").Value = "  -2.05%  "

$ws.Range("D33").Value = "'4.429"
$ws.Range("E33").Value = "  +0.21%  "

$ws.Range("D34").Value = "'0.04930"
$ws.Range("E34").Value = "  -2.66%  "

$ws.Range("D35").Value = "'0.7561"
$ws.Range("E35").Value = "  -1.16%  "

$ws.Range("D37").Value = "'2.734"
$ws.Range("E37").Value = "  +0.30%  "

$ws.Range("D38").Value = "'0.02012"
$ws.Range("E38").Value = "  -0.52%  "

$ws.Range("D39").Value = "'2.703"
$ws.Range("E39").Value = "  -0.18%  "

$ws.Range("D40").Value = "'6.533"
$ws.Range("E40").Value = "  +0.95%  "

$ws.Range("D41").Value = "'78.07"
$ws.Range("E41").Value = "  +12.00%  "

$ws.Range("D42").Value = "'2.119"
$ws.Range("E42").Value = "  -0.40%  "

$ws.Range("D43").Value = "'0.9063"
$ws.Range("E43").Value = "  +3.12%  "

$ws.Range("D44").Value = "'109.48"
$ws.Range("E44").Value = "  +0.14%  "

$ws.Range("D45").Value = "'0.4461"
$ws.Range("E45").Value = "  +0.76%  "

$ws.Range("D46").Value = "'8.163"
$ws.Range("E46").Value = "  +9.57%  "

$ws.Range("E47").Value = "  +0.44%  "

$ws.Range("D48").Value = "1.027.94"
$ws.Range("E48").Value = "  +10.86%  "

$ws.Range("D49").Value = "'9.410"
$ws.Range("E49").Value = "  +0.75%  "

$ws.Range("E50").Value = "  -0.47%  "

$ws.Range("D51").Value = "'35.98"
$ws.Range("E51").Value = "  +0.53%  "
